# Update NATMI LR-pair TPM-derived statistics on the active sheet.
# Only numeric value cells (columns E..T, rows 2..9) change; text columns
# (A..D) and the header row stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, taken from the updated TPM computation.
$updates = @{
    2 = @{ 'E'=3; 'F'=1; 'G'=0.5218183333333334; 'H'=1.565455; 'I'=0.11537035205514; 'J'=0.11537035205514;
           'O'=0.5211722106246275; 'P'=0.5211722106246275; 'Q'=0.2422675545872223; 'R'=2.180407991285;
           'S'=0.06012782142111885; 'T'=0.06012782142111885 }
    3 = @{ 'E'=3; 'F'=1; 'G'=0.5218183333333334; 'H'=1.565455; 'I'=0.11537035205514; 'J'=0.11537035205514;
           'M'=0.426554; 'N'=1.279662; 'O'=0.4788277893753726; 'P'=0.4788277893753726;
           'Q'=0.2225836973566667; 'R'=2.00325327621; 'S'=0.05524253063402117; 'T'=0.05524253063402117 }
    4 = @{ 'I'=0.2874704269761693; 'J'=0.2874704269761693;
           'O'=0.5211722106246275; 'P'=0.5211722106246275; 'Q'=0.6036625191745555; 'R'=5.432962672571;
           'S'=0.1498215979163757; 'T'=0.1498215979163757 }
    5 = @{ 'I'=0.2874704269761693; 'J'=0.2874704269761693;
           'M'=0.426554; 'N'=1.279662; 'O'=0.4788277893753726; 'P'=0.4788277893753726;
           'Q'=0.5546158902806667; 'R'=4.991543012526; 'S'=0.1376488290597936; 'T'=0.1376488290597936 }
    6 = @{ 'G'=2.600134333333333; 'H'=7.800402999999999; 'I'=0.5748713570699702; 'J'=0.5748713570699702;
           'O'=0.5211722106246275; 'P'=0.5211722106246275; 'Q'=1.207179101031222; 'R'=10.864611909281;
           'S'=0.2996069759889359; 'T'=0.2996069759889359 }
    7 = @{ 'G'=2.600134333333333; 'H'=7.800402999999999; 'I'=0.5748713570699702; 'J'=0.5748713570699702;
           'M'=0.426554; 'N'=1.279662; 'O'=0.4788277893753726; 'P'=0.4788277893753726;
           'Q'=1.109097700420667; 'R'=9.981879303786; 'S'=0.2752643810810343; 'T'=0.2752643810810343 }
    8 = @{ 'G'=0.1008076666666667; 'H'=0.302423; 'I'=0.02228786389872057; 'J'=0.02228786389872057;
           'O'=0.5211722106246275; 'P'=0.5211722106246275; 'Q'=0.04680254664677778; 'R'=0.421222919821;
           'S'=0.01161581529819703; 'T'=0.01161581529819703 }
    9 = @{ 'G'=0.1008076666666667; 'H'=0.302423; 'I'=0.02228786389872057; 'J'=0.02228786389872057;
           'M'=0.426554; 'N'=1.279662; 'O'=0.4788277893753726; 'P'=0.4788277893753726;
           'Q'=0.04299991344733334; 'R'=0.386999221026; 'S'=0.01067204860052354; 'T'=0.01067204860052354 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
